$d = $word.ActiveDocument

# Paragraph 20 is the "Requisitos" Heading2 paragraph, paragraph 21 is the
# ListBullet paragraph with the two "Requisito fraco" lines. Both need to be
# removed entirely (including their paragraph marks), leaving the
# Bibliografia paragraph as the last paragraph before the sectPr.
$paraCount = $d.Paragraphs.Count
$startPara = $d.Paragraphs.Item($paraCount - 1)
$rng = $d.Range($startPara.Range.Start, $d.Content.End)
$rng.Delete()
